$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gip"
$ws.Cells.Item(2,3).Value = "Gipr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.01853566666666667
$ws.Cells.Item(2,8).Value = 0.055607
$ws.Cells.Item(2,9).Value = 0.04056456431044909
$ws.Cells.Item(2,10).Value = 0.04056456431044909
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.120294
$ws.Cells.Item(2,14).Value = 0.360882
$ws.Cells.Item(2,15).Value = 0.2714811014452633
$ws.Cells.Item(2,16).Value = 0.2714811014452633
$ws.Cells.Item(2,17).Value = 0.002229729486
$ws.Cells.Item(2,18).Value = 0.020067565374
$ws.Cells.Item(2,19).Value = 0.01101251259864794
$ws.Cells.Item(2,20).Value = 0.01101251259864794

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gip"
$ws.Cells.Item(3,3).Value = "Gipr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.01853566666666667
$ws.Cells.Item(3,8).Value = 0.055607
$ws.Cells.Item(3,9).Value = 0.04056456431044909
$ws.Cells.Item(3,10).Value = 0.04056456431044909
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.2466523333333333
$ws.Cells.Item(3,14).Value = 0.739957
$ws.Cells.Item(3,15).Value = 0.5566482711305432
$ws.Cells.Item(3,16).Value = 0.5566482711305432
$ws.Cells.Item(3,17).Value = 0.004571865433222222
$ws.Cells.Item(3,18).Value = 0.041146788899
$ws.Cells.Item(3,19).Value = 0.02258019459257522
$ws.Cells.Item(3,20).Value = 0.02258019459257522

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Gip"
$ws.Cells.Item(4,3).Value = "Gipr"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.01853566666666667
$ws.Cells.Item(4,8).Value = 0.055607
$ws.Cells.Item(4,9).Value = 0.04056456431044909
$ws.Cells.Item(4,10).Value = 0.04056456431044909
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.07615633333333334
$ws.Cells.Item(4,14).Value = 0.228469
$ws.Cells.Item(4,15).Value = 0.1718706274241937
$ws.Cells.Item(4,16).Value = 0.1718706274241937
$ws.Cells.Item(4,17).Value = 0.001411608409222222
$ws.Cells.Item(4,18).Value = 0.012704475683
$ws.Cells.Item(4,19).Value = 0.006971857119225939
$ws.Cells.Item(4,20).Value = 0.006971857119225939

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gip"
$ws.Cells.Item(5,3).Value = "Gipr"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.4384066666666667
$ws.Cells.Item(5,8).Value = 1.31522
$ws.Cells.Item(5,9).Value = 0.9594354356895509
$ws.Cells.Item(5,10).Value = 0.9594354356895509
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.120294
$ws.Cells.Item(5,14).Value = 0.360882
$ws.Cells.Item(5,15).Value = 0.2714811014452633
$ws.Cells.Item(5,16).Value = 0.2714811014452633
$ws.Cells.Item(5,17).Value = 0.05273769156
$ws.Cells.Item(5,18).Value = 0.47463922404
$ws.Cells.Item(5,19).Value = 0.2604685888466153
$ws.Cells.Item(5,20).Value = 0.2604685888466153

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Gip"
$ws.Cells.Item(6,3).Value = "Gipr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.4384066666666667
$ws.Cells.Item(6,8).Value = 1.31522
$ws.Cells.Item(6,9).Value = 0.9594354356895509
$ws.Cells.Item(6,10).Value = 0.9594354356895509
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.2466523333333333
$ws.Cells.Item(6,14).Value = 0.739957
$ws.Cells.Item(6,15).Value = 0.5566482711305432
$ws.Cells.Item(6,16).Value = 0.5566482711305432
$ws.Cells.Item(6,17).Value = 0.1081340272822222
$ws.Cells.Item(6,18).Value = 0.97320624554
$ws.Cells.Item(6,19).Value = 0.534068076537968
$ws.Cells.Item(6,20).Value = 0.534068076537968

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Gip"
$ws.Cells.Item(7,3).Value = "Gipr"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.4384066666666667
$ws.Cells.Item(7,8).Value = 1.31522
$ws.Cells.Item(7,9).Value = 0.9594354356895509
$ws.Cells.Item(7,10).Value = 0.9594354356895509
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.07615633333333334
$ws.Cells.Item(7,14).Value = 0.228469
$ws.Cells.Item(7,15).Value = 0.1718706274241937
$ws.Cells.Item(7,16).Value = 0.1718706274241937
$ws.Cells.Item(7,17).Value = 0.03338744424222222
$ws.Cells.Item(7,18).Value = 0.30048699818
$ws.Cells.Item(7,19).Value = 0.1648987703049677
$ws.Cells.Item(7,20).Value = 0.1648987703049677
